$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.388.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2936"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.844.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001084"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6791"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.094.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.154"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.411.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.455"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1388"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.373"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.316"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.465"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05616"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.036"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.583"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.232.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.775"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01797"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.459"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9080"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.003.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("E47").Value = "  +3.71%  "
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4014"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.025"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.685"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
